# Remove the "navegador" column (C: navegador/chrome/edge) from Sheet1.
# The "url" column (D) shifts left to become the new column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column C (navegador/chrome/edge), shifting url column left.
$ws.Range("C1:C5").EntireColumn.Delete()

# Update the selected cell to match the recorded view state.
$ws.Range("D13").Select()
